# Fix the mislabeled pseudocode on the "BFS A* - Algorithm" slide: its
# function header was copy/pasted from the Greedy slide and still read
# "Greedy( root , goal)" instead of "AStar( root , goal)".
#
# Commit message: "corrected func name in AStar"
#
# NOTE: the Greedy slide ("BFS Greedy - Algorithm") has a pseudocode box
# with byte-for-byte the same "Greedy( root , goal)" header text (that's
# exactly where it got copy/pasted from), so identifying the right slide
# has to be done by title, not by scanning every slide for that substring.

$p = $ppt.ActivePresentation

function Get-SlideTitle($slide) {
    $title = ""
    try { $title = $slide.Shapes.Title.TextFrame.TextRange.Text } catch {}
    return $title
}

function Is-TargetTitle($title) {
    return ($title.Contains("BFS") -and $title.Contains("A*") -and $title.Contains("Algorithm"))
}

# Primary: slide 13 is "BFS A* - Algorithm" in the original deck.
$targetSlide = $null
if ($p.Slides.Count -ge 13) {
    $candidate = $p.Slides.Item(13)
    if (Is-TargetTitle (Get-SlideTitle $candidate)) {
        $targetSlide = $candidate
    }
}

# Fallback: find the (unique) slide whose title matches "BFS ... A* ...
# Algorithm" in case slide order/count ever changes.
if ($targetSlide -eq $null) {
    for ($si = 1; $si -le $p.Slides.Count; $si++) {
        $slide = $p.Slides.Item($si)
        if (Is-TargetTitle (Get-SlideTitle $slide)) {
            $targetSlide = $slide
            break
        }
    }
}

if ($targetSlide -eq $null) {
    throw "Could not find the 'BFS A* - Algorithm' slide"
}

$targetShape = $null
for ($i = 1; $i -le $targetSlide.Shapes.Count; $i++) {
    $shp = $targetSlide.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $txt = $shp.TextFrame.TextRange.Text
    if ($txt.Contains("Greedy( root , goal)")) {
        $targetShape = $shp
        break
    }
}

if ($targetShape -eq $null) {
    throw "Could not find the pseudocode text box containing 'Greedy( root , goal)'"
}

$tr = $targetShape.TextFrame.TextRange
$full = $tr.Text
$needle = "Greedy( root , goal)"
$idx0 = $full.IndexOf($needle)
if ($idx0 -lt 0) {
    throw "Could not find 'Greedy( root , goal)' text to fix"
}

# Rename the function header, keeping it split the same way the original
# author's edit did: "Greedy( " -> "AStar( ", then "root , goal)" splits
# into "root " + ", goal)".
$funcNameRange = $tr.Characters($idx0 + 1, 8)   # "Greedy( "
$funcNameRange.Text = "AStar( "

$rootRange = $tr.Characters($idx0 + 8, 5)       # "root "
$rootRange.Text = "root "

Write-Output $tr.Text
